# "Excel file reading and writing" -- re-create the First/Last Name
# mini-table (read from columns A:B) over in columns G:I of the
# FirstPage sheet, adding a third ("age") header next to them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstPage")

$ws.Range("G1").Value = "First Name"
$ws.Range("H1").Value = "Last Name"
$ws.Range("I1").Value = "age"

# Auto-size the columns that now hold data (both the original A:D
# block and the new G:I block), the way Excel does right after you
# type into a previously-unused column.
$ws.Columns.Item(1).ColumnWidth = 9.053385416666666
$ws.Columns.Item(2).ColumnWidth = 8.944010416666666
$ws.Columns.Item(3).ColumnWidth = 27.944010416666668
$ws.Columns.Item(4).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(7).ColumnWidth = 9.053385416666666
$ws.Columns.Item(8).ColumnWidth = 8.830729166666666
$ws.Columns.Item(9).ColumnWidth = 3.0533854166666665

# Leave the selection where it ended up after typing the last header.
$ws.Range("J5").Select() | Out-Null
